# Update and delete user information:
#  - "Users to Delete" sheet: refresh the Database ID (UUID) values for the
#    four bulk test users (D2:D5) to reflect the newly (re)created records.
#  - "Summary" sheet: refresh the "Report Generated" timestamp (B6).

$wb = $excel.ActiveWorkbook

$wsUsers = $wb.Worksheets.Item("Users to Delete")
$wsUsers.Range("D2").Value = "0b4b62a1-3e91-45cf-94ac-fb94a84ef9a4"
$wsUsers.Range("D3").Value = "984b1f58-9233-4ec0-a193-30417b5ac9bd"
$wsUsers.Range("D4").Value = "23ed4a6a-2467-48e5-a217-bf69106e9463"
$wsUsers.Range("D5").Value = "e8fc1a8e-bb04-421c-9fda-a519919c51b9"

$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B6").Value = "11/26/2025, 8:36:12 PM"
